$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "otra generacion" tag row now points at the combined
# thermoelectric + interconnection real-time feed instead.
# Using Formula with a leading apostrophe keeps the cell's existing
# "quote prefix" text style (s="3") instead of resetting it.
$ws.Range("E5").Formula = "'/cal/generation_now/termoelectrica + interconnexion"

# Column E needs to be a bit wider to comfortably fit the new, longer tag.
$ws.Columns("E").ColumnWidth = 49.1667

# Move the sheet's remembered selection/active cell.
$ws.Range("E10").Select()
